$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-20: 45185 -> 45204 (2023-10-05)
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
